$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.480.57"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.70%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.704.51"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.54"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5502"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.010"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2744"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.43%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06483"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.15%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07711"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.76%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.695.55"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.560"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5854"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008416"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.99"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.526.39"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.70%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.965"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.93%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.01"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.91"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.28%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.010"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.07"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1333"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +8.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.939"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.91"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06306"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.384"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.88%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.629"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.80%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.29%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.047"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6210"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.28%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.767"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.62%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.123.86"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.168"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8849"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.90%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.854.41"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.92"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.26%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.249"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.20%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.154"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.95%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.12%  "
